# Minor fix in TSP.
# Update the Fitness values (column C) for rows 2-12 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(3565, 3565, 3565, 3565, 3565, 3565, 3617, 3799, 3799, 3851, 3851)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
